$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("C18").Value = "'0.5117"
$ws.Range("D18").Value = "'0.206"
$ws.Range("E18").Value = "'1.268"
$ws.Range("C25").Value = "'0.8986"
$ws.Range("D25").Value = "'0.782"
$ws.Range("E25").Value = "'1.032"
$ws.Range("C37").Value = "'0.6229"
$ws.Range("D37").Value = "'0.415"
$ws.Range("E37").Value = "'0.935"
$ws.Range("C38").Value = "'0.7986"
$ws.Range("E38").Value = "'7.986"
$ws.Range("C45").Value = "'0.8146"
$ws.Range("D45").Value = "'0.645"
$ws.Range("E45").Value = "'1.029"
$ws.Range("C50").Value = "'1.0185"
$ws.Range("D50").Value = "'0.923"
$ws.Range("E50").Value = "'1.124"
$ws.Range("C52").Value = "'0.9801"
$ws.Range("D52").Value = "'0.779"
$ws.Range("E52").Value = "'1.233"
$ws.Range("C54").Value = "'0.9337"
$ws.Range("D54").Value = "'0.610"
$ws.Range("E54").Value = "'1.430"
$ws.Range("C56").Value = "'0.8419"
$ws.Range("D56").Value = "'0.374"
$ws.Range("E56").Value = "'1.898"
$ws.Range("C57").Value = "'0.2839"
$ws.Range("D57").Value = "'0.108"
$ws.Range("E57").Value = "'0.747"
$ws = $wb.Worksheets.Item(2)
$ws.Range("C5").Value = "'0.9706"
$ws.Range("D5").Value = "'0.748"
$ws.Range("E5").Value = "'1.260"
$ws.Range("C6").Value = "'0.7755"
$ws.Range("D6").Value = "'0.692"
$ws.Range("E6").Value = "'0.869"
$ws.Range("C7").Value = "'0.8445"
$ws.Range("D7").Value = "'0.779"
$ws.Range("E7").Value = "'0.916"
$ws.Range("C14").Value = "'0.9388"
$ws.Range("D14").Value = "'0.821"
$ws.Range("E14").Value = "'1.074"
$ws.Range("C17").Value = "'0.6768"
$ws.Range("D17").Value = "'0.517"
$ws.Range("E17").Value = "'0.887"
$ws.Range("C18").Value = "'1.0291"
$ws.Range("D18").Value = "'0.644"
$ws.Range("E18").Value = "'1.644"
$ws.Range("C20").Value = "'0.7909"
$ws.Range("E20").Value = "'7.909"
$ws.Range("C21").Value = "'0.2015"
$ws.Range("D21").Value = "'0.112"
$ws.Range("E21").Value = "'0.364"
$ws.Range("C25").Value = "'0.8326"
$ws.Range("D25").Value = "'0.732"
$ws.Range("E25").Value = "'0.947"
$ws.Range("C26").Value = "'0.8584"
$ws.Range("D26").Value = "'0.802"
$ws.Range("E26").Value = "'0.919"
$ws.Range("C27").Value = "'0.9624"
$ws.Range("D27").Value = "'0.915"
$ws.Range("E27").Value = "'1.013"
$ws.Range("C35").Value = "'1.0084"
$ws.Range("D35").Value = "'0.837"
$ws.Range("E35").Value = "'1.214"
$ws.Range("C37").Value = "'0.8990"
$ws.Range("D37").Value = "'0.652"
$ws.Range("E37").Value = "'1.239"
$ws.Range("C40").Value = "'0.7383"
$ws.Range("D40").Value = "'0.379"
$ws.Range("E40").Value = "'1.436"
$ws.Range("C41").Value = "'0.5681"
$ws.Range("D41").Value = "'0.288"
$ws.Range("E41").Value = "'1.122"
$ws.Range("C45").Value = "'0.8578"
$ws.Range("D45").Value = "'0.663"
$ws.Range("E45").Value = "'1.110"
$ws.Range("C55").Value = "'0.8265"
$ws.Range("D55").Value = "'0.663"
$ws.Range("E55").Value = "'1.031"
$ws.Range("C59").Value = "'0.7184"
$ws.Range("D59").Value = "'0.246"
$ws.Range("E59").Value = "'2.098"
$ws = $wb.Worksheets.Item(3)
$ws.Range("C5").Value = "'0.6596"
$ws.Range("D5").Value = "'0.453"
$ws.Range("E5").Value = "'0.961"
$ws.Range("C9").Value = "'0.8578"
$ws.Range("D9").Value = "'0.777"
$ws.Range("E9").Value = "'0.947"
$ws.Range("C10").Value = "'0.9434"
$ws.Range("D10").Value = "'0.846"
$ws.Range("E10").Value = "'1.052"
$ws.Range("C15").Value = "'0.8445"
$ws.Range("D15").Value = "'0.663"
$ws.Range("E15").Value = "'1.075"
$ws.Range("C16").Value = "'0.9015"
$ws.Range("D16").Value = "'0.701"
$ws.Range("E16").Value = "'1.159"
$ws.Range("C18").Value = "'0.8994"
$ws.Range("D18").Value = "'0.560"
$ws.Range("E18").Value = "'1.444"
$ws.Range("C20").Value = "'0.6704"
$ws.Range("D20").Value = "'0.326"
$ws.Range("E20").Value = "'1.379"
$ws.Range("C25").Value = "'0.7706"
$ws.Range("D25").Value = "'0.643"
$ws.Range("E25").Value = "'0.923"
$ws.Range("C28").Value = "'1.0170"
$ws.Range("D28").Value = "'0.959"
$ws.Range("E28").Value = "'1.078"
$ws.Range("C29").Value = "'0.9822"
$ws.Range("D29").Value = "'0.929"
$ws.Range("E29").Value = "'1.038"
$ws.Range("C34").Value = "'0.9641"
$ws.Range("D34").Value = "'0.853"
$ws.Range("E34").Value = "'1.090"
$ws.Range("C35").Value = "'0.8534"
$ws.Range("D35").Value = "'0.744"
$ws.Range("E35").Value = "'0.979"
$ws.Range("C38").Value = "'1.0533"
$ws.Range("D38").Value = "'0.760"
$ws.Range("E38").Value = "'1.460"
$ws.Range("C46").Value = "'0.7044"
$ws.Range("D46").Value = "'0.564"
$ws.Range("E46").Value = "'0.880"
$ws.Range("C48").Value = "'0.7884"
$ws.Range("D48").Value = "'0.692"
$ws.Range("E48").Value = "'0.898"
$ws.Range("C53").Value = "'0.9340"
$ws.Range("D53").Value = "'0.777"
$ws.Range("E53").Value = "'1.123"
$ws.Range("C54").Value = "'1.0625"
$ws.Range("D54").Value = "'0.874"
$ws.Range("E54").Value = "'1.292"
$ws.Range("C59").Value = "'0.6908"
$ws.Range("D59").Value = "'0.347"
$ws.Range("E59").Value = "'1.375"
$ws.Range("C78").Value = "'0.9315"
$ws.Range("D78").Value = "'0.663"
$ws.Range("E78").Value = "'1.308"
$ws.Range("C98").Value = "'0.8844"
$ws.Range("D98").Value = "'0.589"
$ws.Range("E98").Value = "'1.327"
$ws.Range("C101").Value = "'0.5380"
$ws.Range("D101").Value = "'0.264"
$ws.Range("E101").Value = "'1.095"
$ws = $wb.Worksheets.Item(4)
$ws.Range("C19").Value = "'0.5340"
$ws.Range("D19").Value = "'0.275"
$ws.Range("E19").Value = "'1.038"
$ws.Range("C20").Value = "'0.8316"
$ws.Range("D20").Value = "'0.325"
$ws.Range("E20").Value = "'2.126"
$ws.Range("C21").Value = "'0.4408"
$ws.Range("D21").Value = "'0.168"
$ws.Range("E21").Value = "'1.159"
$ws.Range("C25").Value = "'1.1495"
$ws.Range("D25").Value = "'0.468"
$ws.Range("E25").Value = "'2.821"
$ws.Range("C26").Value = "'0.9550"
$ws.Range("D26").Value = "'0.735"
$ws.Range("E26").Value = "'1.240"
$ws.Range("C34").Value = "'0.4953"
$ws.Range("D34").Value = "'0.357"
$ws.Range("E34").Value = "'0.687"
$ws.Range("C37").Value = "'0.5117"
$ws.Range("D37").Value = "'0.249"
$ws.Range("E37").Value = "'1.051"
$ws.Range("C39").Value = "'0.6231"
$ws.Range("D39").Value = "'0.289"
$ws.Range("E39").Value = "'1.344"
$ws.Range("C46").Value = "'0.9561"
$ws.Range("D46").Value = "'0.847"
$ws.Range("E46").Value = "'1.079"
$ws.Range("C59").Value = "'0.5293"
$ws.Range("D59").Value = "'0.270"
$ws.Range("E59").Value = "'1.037"
$ws.Range("C61").Value = "'0.3897"
$ws.Range("D61").Value = "'0.153"
$ws.Range("E61").Value = "'0.992"
$ws.Range("C74").Value = "'0.9795"
$ws.Range("D74").Value = "'0.768"
$ws.Range("E74").Value = "'1.249"
$ws.Range("C79").Value = "'0.5113"
$ws.Range("D79").Value = "'0.260"
$ws.Range("E79").Value = "'1.005"
$ws.Range("C80").Value = "'0.6230"
$ws.Range("D80").Value = "'0.234"
$ws.Range("E80").Value = "'1.658"
$ws.Range("C81").Value = "'0.2063"
$ws.Range("D81").Value = "'0.075"
$ws.Range("E81").Value = "'0.565"
